$wbx = $excel.ActiveWorkbook
$ws = $wbx.ActiveSheet

# Insert a new row at 469, shifting existing rows 469:570 down to 470:571.
$ws.Rows.Item(469).Insert()

# Populate the newly inserted row 469 with the new weekly data entry.
$ws.Range("A469").Value = 7
$ws.Range("B469").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C469").Value = "Ñuble"
$ws.Range("D469").Value = 45275
$ws.Range("E469").Value = 16
$ws.Range("F469").Value = 100112006
$ws.Range("G469").Value = "Repollo"
$ws.Range("H469").Value = "Morada(o)"
$ws.Range("I469").Value = "Primera"
$ws.Range("J469").Value = 200
$ws.Range("K469").Value = 1500
$ws.Range("L469").Value = 1500
$ws.Range("M469").Value = 1500
$ws.Range("N469").Value = "$/unidad"
$ws.Range("O469").Value = "Región del Maule"
$ws.Range("P469").Value = 1500
$ws.Range("Q469").Value = 1
$ws.Range("R469").Value = "Hortaliza"
